$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1137.7858
$ws.Range("I28").Value = 987.9524
$ws.Range("J28").Value = 1587.2858
$ws.Range("K28").Value = 987.9524
$ws.Range("L28").Value = 1587.2858
$ws.Range("M28").Value = -502.9524
$ws.Range("N28").Value = -2557.2858
$ws.Range("H33").Value = 126353.25
$ws.Range("I33").Value = 367.3
$ws.Range("K33").Value = 367.3
$ws.Range("M33").Value = -138.3
$ws.Range("H116").Value = 6129.1665
$ws.Range("I116").Value = 3700.625
$ws.Range("K116").Value = 3700.625
$ws.Range("M116").Value = -258.625
$ws.Range("H137").Value = 2046.238
$ws.Range("I137").Value = 1591.1538
$ws.Range("J137").Value = 2785.75
$ws.Range("K137").Value = 4773.4614
$ws.Range("L137").Value = 8357.25
$ws.Range("M137").Value = -2223.4614
$ws.Range("N137").Value = -13457.25
$ws.Range("H141").Value = 962.3333
$ws.Range("I141").Value = 943.5
$ws.Range("K141").Value = 2830.5
$ws.Range("M141").Value = 2349.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6468.8057
$ws.Range("I61").Value = 4811.1924
$ws.Range("J61").Value = 10778.6
$ws.Range("K61").Value = 4811.1924
$ws.Range("L61").Value = 10778.6
$ws.Range("M61").Value = -4599.1924
$ws.Range("N61").Value = -11202.6
$ws.Range("H74").Value = 2056.4314
$ws.Range("I74").Value = 2015.898
$ws.Range("K74").Value = 2015.898
$ws.Range("M74").Value = -1141.898
$ws.Range("H77").Value = 2056.4314
$ws.Range("I77").Value = 2015.898
$ws.Range("K77").Value = 10079.49
$ws.Range("M77").Value = -5711.49
$ws.Range("H122").Value = 2884.0715
$ws.Range("I122").Value = 2924.2222
$ws.Range("K122").Value = 8772.6666
$ws.Range("M122").Value = -6322.6666
$ws.Range("H132").Value = 8428.944
$ws.Range("I132").Value = 8370.1875
$ws.Range("K132").Value = 25110.5625
$ws.Range("M132").Value = -22580.5625
$ws.Range("H136").Value = 6468.8057
$ws.Range("I136").Value = 4811.1924
$ws.Range("J136").Value = 10778.6
$ws.Range("K136").Value = 14433.5772
$ws.Range("L136").Value = 32335.8
$ws.Range("M136").Value = -11883.5772
$ws.Range("N136").Value = -37435.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 660
$ws.Range("I22").Value = 660
$ws.Range("K22").Value = 660
$ws.Range("M22").Value = -487
$ws.Range("H105").Value = 2979.111
$ws.Range("I105").Value = 2101.5
$ws.Range("K105").Value = 2101.5
$ws.Range("M105").Value = -354.5
$ws.Range("H107").Value = 2547.4666
$ws.Range("I107").Value = 1692.9445
$ws.Range("J107").Value = 3829.25
$ws.Range("K107").Value = 1692.9445
$ws.Range("L107").Value = 3829.25
$ws.Range("M107").Value = 227.0554999999999
$ws.Range("N107").Value = -7669.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5427.423
$ws.Range("I31").Value = 6268.375
$ws.Range("J31").Value = 4081.9
$ws.Range("K31").Value = 6268.375
$ws.Range("L31").Value = 4081.9
$ws.Range("M31").Value = -5973.375
$ws.Range("N31").Value = -4671.9
$ws.Range("H34").Value = 5427.423
$ws.Range("I34").Value = 6268.375
$ws.Range("J34").Value = 4081.9
$ws.Range("K34").Value = 6268.375
$ws.Range("L34").Value = 4081.9
$ws.Range("M34").Value = -6066.375
$ws.Range("N34").Value = -4485.9
$ws.Range("H99").Value = 7540.0835
$ws.Range("I99").Value = 7157.1763
$ws.Range("J99").Value = 8470
$ws.Range("K99").Value = 7157.1763
$ws.Range("L99").Value = 8470
$ws.Range("M99").Value = -5659.1763
$ws.Range("N99").Value = -11466
$ws.Range("H107").Value = 272.82608
$ws.Range("I107").Value = 259.89474
$ws.Range("J107").Value = 334.25
$ws.Range("K107").Value = 259.89474
$ws.Range("L107").Value = 334.25
$ws.Range("M107").Value = 1660.10526
$ws.Range("N107").Value = -4174.25
$ws.Range("H126").Value = 7540.0835
$ws.Range("I126").Value = 7157.1763
$ws.Range("J126").Value = 8470
$ws.Range("K126").Value = 21471.5289
$ws.Range("L126").Value = 25410
$ws.Range("M126").Value = -19001.5289
$ws.Range("N126").Value = -30350
$ws.Range("H134").Value = 1787.2916
$ws.Range("I134").Value = 852.45
$ws.Range("K134").Value = 2557.35
$ws.Range("M134").Value = -22.35000000000036

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 603.3871
$ws.Range("I2").Value = 1090.4667
$ws.Range("K2").Value = 6542.8002
$ws.Range("M2").Value = -6429.8002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 7627.25
$ws.Range("I22").Value = 11250
$ws.Range("J22").Value = 4004.5
$ws.Range("K22").Value = 11250
$ws.Range("L22").Value = 4004.5
$ws.Range("M22").Value = -10721
$ws.Range("N22").Value = -5062.5
$ws.Range("H97").Value = 675.46155
$ws.Range("J97").Value = 779.2857
$ws.Range("L97").Value = 779.2857
$ws.Range("N97").Value = -1771.2857
$ws.Range("H113").Value = 6651.5713
$ws.Range("I113").Value = 7108.909
$ws.Range("K113").Value = 7108.909
$ws.Range("M113").Value = -4938.909
$ws.Range("H132").Value = 4739.909
$ws.Range("I132").Value = 4758.263
$ws.Range("K132").Value = 14274.789
$ws.Range("M132").Value = -11744.789

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 30088.762
$ws.Range("I7").Value = 31992
$ws.Range("J7").Value = 23998.4
$ws.Range("K7").Value = 31992
$ws.Range("L7").Value = 23998.4
$ws.Range("M7").Value = -31880
$ws.Range("N7").Value = -24222.4
$ws.Range("H22").Value = 3232.6428
$ws.Range("I22").Value = 2876.5
$ws.Range("J22").Value = 3707.5
$ws.Range("K22").Value = 2876.5
$ws.Range("L22").Value = 3707.5
$ws.Range("M22").Value = -2581.5
$ws.Range("N22").Value = -4297.5
$ws.Range("H25").Value = 8066.3335
$ws.Range("I25").Value = 5999.5
$ws.Range("K25").Value = 5999.5
$ws.Range("M25").Value = -5769.5
$ws.Range("H27").Value = 3232.6428
$ws.Range("I27").Value = 2876.5
$ws.Range("J27").Value = 3707.5
$ws.Range("K27").Value = 2876.5
$ws.Range("L27").Value = 3707.5
$ws.Range("M27").Value = -2769.5
$ws.Range("N27").Value = -3921.5
$ws.Range("H40").Value = 4017.2917
$ws.Range("I40").Value = 3444.5386
$ws.Range("K40").Value = 3444.5386
$ws.Range("M40").Value = -3308.5386
$ws.Range("H93").Value = 1262.25
$ws.Range("I93").Value = 1262.25
$ws.Range("K93").Value = 1262.25
$ws.Range("M93").Value = -14.25
$ws.Range("H126").Value = 30088.762
$ws.Range("I126").Value = 31992
$ws.Range("J126").Value = 23998.4
$ws.Range("K126").Value = 95976
$ws.Range("L126").Value = 71995.20000000001
$ws.Range("M126").Value = -93506
$ws.Range("N126").Value = -76935.20000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3364.3572
$ws.Range("I126").Value = 2683.4167
$ws.Range("J126").Value = 7450
$ws.Range("K126").Value = 8050.250100000001
$ws.Range("L126").Value = 22350
$ws.Range("M126").Value = -5580.250100000001
$ws.Range("N126").Value = -27290
$ws.Range("H132").Value = 4786.375
$ws.Range("I132").Value = 4861.25
$ws.Range("K132").Value = 14583.75
$ws.Range("M132").Value = -12053.75

Write-Output "Applied 184 cell updates"